$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 49
$ws.Range("F5").Value = 1064
$ws.Range("F7").Value = 611
$ws.Range("F8").Value = 576
$ws.Range("F9").Value = 1503
$ws.Range("F11").Value = 1398
$ws.Range("F12").Value = 3042
$ws.Range("F13").Value = 542
$ws.Range("F14").Value = 1701
$ws.Range("F15").Value = 1710
$ws.Range("F17").Value = 255
$ws.Range("F18").Value = 1430
$ws.Range("F21").Value = 1162
$ws.Range("F22").Value = 378
$ws.Range("F23").Value = 420
$ws.Range("F24").Value = 46
$ws.Range("F25").Value = 3748
$ws.Range("F26").Value = 717
$ws.Range("F27").Value = 565
$ws.Range("F28").Value = 1597
$ws.Range("F29").Value = 1
$ws.Range("F30").Value = 65

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 15
$ws.Range("F9").Value = 38
$ws.Range("F13").Value = 91

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 25

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 49
$ws.Range("F4").Value = 25
$ws.Range("F12").Value = 15
$ws.Range("F14").Value = 38
$ws.Range("F16").Value = 1064
$ws.Range("F18").Value = 611
$ws.Range("F19").Value = 576
$ws.Range("F20").Value = 1503
$ws.Range("F22").Value = 1398
$ws.Range("F23").Value = 3042
$ws.Range("F24").Value = 542
$ws.Range("F25").Value = 1701
$ws.Range("F26").Value = 1711
$ws.Range("F28").Value = 255
$ws.Range("F29").Value = 1430
$ws.Range("F34").Value = 1162
$ws.Range("F35").Value = 378
$ws.Range("F36").Value = 420
$ws.Range("F37").Value = 46
$ws.Range("F38").Value = 3749
$ws.Range("F39").Value = 717
$ws.Range("F40").Value = 565
$ws.Range("F41").Value = 1597
$ws.Range("F42").Value = 91
$ws.Range("F44").Value = 1
$ws.Range("F45").Value = 65
